$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (D) / Volume(1h) (E) values scraped for this run, keyed by row number.
$updates = @{
    "2" = @{ D = "26.629.26"; E = "  +1.03%  " }
    "3" = @{ D = "1.826.65"; E = "  +1.76%  " }
    "4" = @{ D = "1.008"; E = "  +0.26%  " }
    "5" = @{ D = "1.007"; E = "  +0.28%  " }
    "6" = @{ D = "308.54"; E = "  +0.47%  " }
    "7" = @{ D = "0.4675"; E = "  +3.78%  " }
    "8" = @{ E = "  +0.15%  " }
    "9" = @{ D = "0.07133"; E = "  +0.62%  " }
    "10" = @{ D = "0.9020"; E = "  +1.96%  " }
    "11" = @{ D = "0.07755"; E = "  +0.25%  " }
    "12" = @{ D = "19.42"; E = "  +0.09%  " }
    "13" = @{ D = "1.819.38"; E = "  +1.46%  " }
    "14" = @{ D = "5.265"; E = "  -0.32%  " }
    "15" = @{ D = "6.345"; E = "  +0.29%  " }
    "16" = @{ D = "87.52"; E = "  +3.11%  " }
    "17" = @{ D = "1.010"; E = "  +0.39%  " }
    "18" = @{ D = "0.000008552"; E = "  +0.43%  " }
    "19" = @{ D = "1.007"; E = "  +0.29%  " }
    "20" = @{ D = "26.657.56"; E = "  +1.07%  " }
    "21" = @{ D = "14.21"; E = "  -0.37%  " }
    "22" = @{ D = "5.020"; E = "  +0.81%  " }
    "23" = @{ D = "10.55"; E = "  +0.21%  " }
    "24" = @{ D = "1.911"; E = "  -2.85%  " }
    "25" = @{ D = "152.97"; E = "  +1.01%  " }
    "26" = @{ D = "17.93"; E = "  +0.51%  " }
    "27" = @{ D = "1.972"; E = "  -2.49%  " }
    "28" = @{ D = "113.87"; E = "  +1.71%  " }
    "29" = @{ D = "4.871"; E = "  -0.33%  " }
    "30" = @{ D = "0.08811"; E = "  +1.56%  " }
    "31" = @{ D = "3.146"; E = "  +2.88%  " }
    "32" = @{ D = "2.817"; E = "  +2.67%  " }
    "33" = @{ D = "1.160"; E = "  +4.85%  " }
    "34" = @{ D = "0.7353"; E = "  +1.42%  " }
    "35" = @{ D = "4.444"; E = "  -0.03%  " }
    "36" = @{ D = "1.081"; E = "  +1.32%  " }
    "37" = @{ D = "0.01929"; E = "  -0.09%  " }
    "38" = @{ E = "  +1.29%  " }
    "39" = @{ D = "2.913"; E = "  +1.84%  " }
    "40" = @{ D = "6.887"; E = "  +0.30%  " }
    "41" = @{ D = "0.5058"; E = "  -0.25%  " }
    "42" = @{ D = "0.1496"; E = "  -1.06%  " }
    "43" = @{ D = "8.033"; E = "  +0.33%  " }
    "44" = @{ D = "1.008"; E = "  +0.34%  " }
    "45" = @{ D = "0.4665"; E = "  +0.64%  " }
    "46" = @{ D = "10.01"; E = "  +1.53%  " }
    "47" = @{ D = "98.00"; E = "  -3.03%  " }
    "48" = @{ D = "1.568"; E = "  -0.86%  " }
    "49" = @{ D = "0.06051"; E = "  +1.55%  " }
    "50" = @{ D = "63.95"; E = "  -0.27%  " }
    "51" = @{ D = "35.84"; E = "  -0.36%  " }
}

# The Price column holds plain numeric-looking text (e.g. "1.008"); force the cells
# to Text format first so Excel does not silently coerce the assigned strings into
# numbers, then strip the temporary format back off so styling is unaffected.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $ws.Range("D$row").Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
}

$priceRange.ClearFormats()

